# Sprint(46) test case report updates - crabvpn test and jenkins run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day 3 (rows 14-17): execution/review totals corrected 3120 -> 3160
$ws.Range("C16").Value = 3160
$ws.Range("C17").Value = 3160

# Day 4 (rows 20-23): totals filled in (were previously blank)
$ws.Range("C21").Value = 7300
$ws.Range("C22").Value = 3190
$ws.Range("C23").Value = 3190

# Day 6 (rows 26-29): Total Execution was casual leave that day
$ws.Range("C28").Value = "casual leave"

# Day 7 (rows 32-35): Total Execution was casual leave that day
$ws.Range("C34").Value = "casual leave"

# Day 8 (rows 38-41): totals filled in (were previously blank)
$ws.Range("C39").Value = 7300
$ws.Range("C40").Value = 3190
$ws.Range("C41").Value = 3190

# Day 9 (rows 44-47): totals filled in (were previously blank)
$ws.Range("C45").Value = 7300
$ws.Range("C46").Value = 3220
$ws.Range("C47").Value = 3220

# Reflect the scrolled position/selection that was active when the
# workbook was saved (view scrolled down to row 32, C47 selected)
$win = $excel.ActiveWindow
$win.ScrollRow = 32
$win.ScrollColumn = 1
$ws.Range("C47").Select()
